$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New consolidated values for rows 2-9
$values = @(
    "('Avatar', ['Token Creature — Avatar', 'This creature’s power and toughness are each equal to your life total.', '*/*'])",
    "('Beast', ['Token Creature — Beast', '3/3'])",
    "('Gargoyle', ['Token Artifact Creature — Gargoyle', 'Flying', '3/4'])",
    "('Goblin', ['Token Creature — Goblin', '1/1'])",
    "('Insect', ['Token Creature — Insect', '1/1'])",
    "('Soldier', ['Token Creature — Soldier', '1/1'])",
    "('Wolf', ['Token Creature — Wolf', '2/2'])",
    "('Zombie', ['Token Creature — Zombie', '2/2'])"
)

# Write the new consolidated values into A2:A9
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Clear out the old rows 10-27 which no longer have data
$ws.Range("A10:A27").Clear()
